$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.950.19'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '1.764.60'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '328.40'
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '0.4656'
$ws.Range("E7").Value = '  +0.67%  '
$ws.Range("D8").Value = '0.3509'
$ws.Range("E8").Value = '  -2.21%  '
$ws.Range("D9").Value = '43.25'
$ws.Range("E9").Value = '  +3.55%  '
$ws.Range("D10").Value = '0.07354'
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("D11").Value = '1.082'
$ws.Range("E11").Value = '  -1.69%  '
$ws.Range("D12").Value = '0.9999'
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Value = '20.63'
$ws.Range("E13").Value = '  -0.95%  '
$ws.Range("D14").Value = '5.989'
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").Value = '7.154'
$ws.Range("E15").Value = '  -1.31%  '
$ws.Range("D16").Value = '1.762.65'
$ws.Range("E16").Value = '  -0.87%  '
$ws.Range("D17").Value = '92.44'
$ws.Range("E17").Value = '  -1.31%  '
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("E20").Value = '  -0.07%  '
$ws.Range("D21").Value = '16.85'
$ws.Range("E21").Value = '  -1.38%  '
$ws.Range("D22").Value = '5.763'
$ws.Range("E22").Value = '  -0.27%  '
$ws.Range("D23").Value = '27.975.80'
$ws.Range("D24").Value = '11.13'
$ws.Range("E24").Value = '  -1.41%  '
$ws.Range("D25").Value = '2.152'
$ws.Range("E25").Value = '  +3.46%  '
$ws.Range("D26").Value = '162.23'
$ws.Range("E26").Value = '  -1.40%  '
$ws.Range("D27").Value = '20.02'
$ws.Range("E27").Value = '  -1.51%  '
$ws.Range("D28").Value = '1.965.61'
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("D29").Value = '2.164'
$ws.Range("E29").Value = '  -0.30%  '
$ws.Range("D30").Value = '122.76'
$ws.Range("E30").Value = '  -2.58%  '
$ws.Range("D31").Value = '1.069'
$ws.Range("E31").Value = '  -2.32%  '
$ws.Range("D32").Value = '0.09275'
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("D33").Value = '3.643'
$ws.Range("E33").Value = '  -0.84%  '
$ws.Range("D34").Value = '5.554'
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("D35").Value = '11.67'
$ws.Range("E35").Value = '  -1.14%  '
$ws.Range("E36").Value = '  -1.08%  '
$ws.Range("D37").Value = '0.06066'
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("D38").Value = '0.2061'
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("D39").Value = '4.912'
$ws.Range("E39").Value = '  -1.01%  '
$ws.Range("D40").Value = '0.6117'
$ws.Range("E40").Value = '  -3.00%  '
$ws.Range("D41").Value = '1.181'
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '1.372'
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("D43").Value = '7.777'
$ws.Range("E43").Value = '  +0.00%  '
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("D45").Value = '3.736'
$ws.Range("E45").Value = '  +0.11%  '
$ws.Range("D46").Value = '0.5788'
$ws.Range("E46").Value = '  -1.61%  '
$ws.Range("D47").Value = '122.89'
$ws.Range("E47").Value = '  +0.57%  '
$ws.Range("D48").Value = '1.923'
$ws.Range("E48").Value = '  -1.23%  '
$ws.Range("D49").Value = '0.06813'
$ws.Range("E49").Value = '  -1.74%  '
$ws.Range("D50").Value = '1.122'
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("D51").Value = '72.02'
$ws.Range("E51").Value = '  -0.25%  '
